$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Student info
$ws.Range("C4").Value = "pepster"
$ws.Range("C5").Value = "Plamena Nedelcheva Georgieva"

# GitHub profile link -> becomes a real hyperlink
$ws.Hyperlinks.Add($ws.Range("C7"), "https://github.com/plmng/IssueTrackingSystem") | Out-Null

# Scores filled in
$ws.Range("C23").Value = 10
$ws.Range("C28").Value = 10
$ws.Range("E28").Value = "not set as buttons but links in navigations"

# Leave the selection where the author left it when saving
$ws.Range("E29").Select() | Out-Null
